$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + report date range) ---
$ws.Range("A8").Value = "Volume 29   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/24/2022  Through  10/30/2022"

# --- Weekly crime statistics table refresh (rows 14-30) ---
# Plain value updates (no type/style change needed)
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -66.666666666666
$ws.Range("I14").Value = 8
$ws.Range("K14").Value = -33.333333333333
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = -42.857142857142
$ws.Range("N14").Value = -60
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 30
$ws.Range("K15").Value = -23.076923076923
$ws.Range("L15").Value = -34.782608695652
$ws.Range("M15").Value = -37.5
$ws.Range("N15").Value = -55.882352941176
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 31.25
$ws.Range("I16").Value = 222
$ws.Range("J16").Value = 166
$ws.Range("K16").Value = 33.734939759036
$ws.Range("L16").Value = 12.690355329949
$ws.Range("M16").Value = -36.571428571428
$ws.Range("N16").Value = -80.457746478873
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 64
$ws.Range("G17").Value = 44
$ws.Range("H17").Value = 45.454545454545
$ws.Range("I17").Value = 656
$ws.Range("J17").Value = 498
$ws.Range("K17").Value = 31.726907630522
$ws.Range("L17").Value = 29.133858267716
$ws.Range("M17").Value = 60.391198044009
$ws.Range("N17").Value = -36.372453928225
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 46.666666666666
$ws.Range("I18").Value = 260
$ws.Range("J18").Value = 207
$ws.Range("K18").Value = 25.603864734299
$ws.Range("L18").Value = -6.810035842293
$ws.Range("M18").Value = -49.709864603481
$ws.Range("N18").Value = -90.931287059644
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100
$ws.Range("G19").Value = 104
$ws.Range("H19").Value = -3.846153846153
$ws.Range("I19").Value = 1121
$ws.Range("J19").Value = 891
$ws.Range("K19").Value = 25.813692480359
$ws.Range("L19").Value = 41.540404040404
$ws.Range("M19").Value = 38.737623762376
$ws.Range("N19").Value = -18.591140159767
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 280
$ws.Range("F20").Value = 48
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 391
$ws.Range("J20").Value = 208
$ws.Range("K20").Value = 87.980769230769
$ws.Range("L20").Value = 98.477157360406
$ws.Range("M20").Value = 31.208053691275
$ws.Range("N20").Value = -90.290538862676
$ws.Range("C21").Value = 68
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 51.111111111111
$ws.Range("F21").Value = 260
$ws.Range("G21").Value = 213
$ws.Range("H21").Value = 22.06572769953
$ws.Range("I21").Value = 2688
$ws.Range("J21").Value = 2021
$ws.Range("K21").Value = 33.003463631865
$ws.Range("L21").Value = 32.088452088452
$ws.Range("M21").Value = 9.983633387888
$ws.Range("N21").Value = -74.463233897016
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -75
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -55.555555555555
$ws.Range("I23").Value = 81
$ws.Range("J23").Value = 89
$ws.Range("K23").Value = -8.988764044943
$ws.Range("L23").Value = 28.571428571428
$ws.Range("M23").Value = 39.655172413793
$ws.Range("C24").Value = 81
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = 92.857142857142
$ws.Range("F24").Value = 349
$ws.Range("G24").Value = 208
$ws.Range("H24").Value = 67.788461538461
$ws.Range("I24").Value = 3260
$ws.Range("J24").Value = 2049
$ws.Range("K24").Value = 59.102000976085
$ws.Range("L24").Value = 52.478952291861
$ws.Range("M24").Value = 2.969046114971
$ws.Range("C25").Value = 27
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -3.571428571428
$ws.Range("F25").Value = 112
$ws.Range("G25").Value = 101
$ws.Range("H25").Value = 10.89108910891
$ws.Range("I25").Value = 1378
$ws.Range("J25").Value = 1119
$ws.Range("K25").Value = 23.145665773011
$ws.Range("L25").Value = 32.5
$ws.Range("M25").Value = -19.320843091334
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 60
$ws.Range("I26").Value = 63
$ws.Range("J26").Value = 76
$ws.Range("K26").Value = -17.105263157894
$ws.Range("L26").Value = -16
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 12
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 144
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = 58.241758241758
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 32
$ws.Range("K28").Value = -3.030303030303
$ws.Range("L28").Value = -13.513513513513
$ws.Range("M28").Value = 23.076923076923
$ws.Range("N28").Value = -67.01030927835
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 26
$ws.Range("K29").Value = -13.333333333333
$ws.Range("L29").Value = -21.212121212121
$ws.Range("M29").Value = 8.333333333333
$ws.Range("N29").Value = -69.411764705882

# Cells that become placeholder text ("0" / "***.*") - set value then fix style
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "***.*"

# Cells that become real numbers (were placeholder text) - set value then fix style
$ws.Range("C15").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50

# --- Fix number formats / styles for cells that changed between text and numeric ---
$countFormatRef = $ws.Range("F16")   # style s=15 (#,##0) used by count columns
$pctFormatRef   = $ws.Range("H16")   # style s=16 (#,##0.0) used by percent columns
$textFormatRef  = $ws.Range("C22")   # style s=14 (General) used by placeholder text cells

$countFormatRef.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D26").PasteSpecial(-4122)

$pctFormatRef.Copy()
$ws.Range("E26").PasteSpecial(-4122)

$textFormatRef.Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row label fix: A30 label unaffected; A15/A37 keep "Rape" text (no value change) ---
